# A new daily price record for Frambuesa (raspberry) was reported and needs
# to be inserted as a new row right before the current row 98, pushing all
# subsequent rows down by one (rows 98-175 become rows 99-176).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 98 (existing row 98 and everything
# below shifts down by one row).
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new record's data.
$ws.Range("A98").Value = 9
$ws.Range("B98").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C98").Value = "Metropolitana"
$ws.Range("D98").Value = 45096
$ws.Range("E98").Value = 13
$ws.Range("F98").Value = "Fruta"
$ws.Range("G98").Value = 100101
$ws.Range("H98").Value = "Berries"
$ws.Range("I98").Value = 100101004
$ws.Range("J98").Value = "Frambuesa"
$ws.Range("K98").Value = "Sin especificar"
$ws.Range("L98").Value = "Primera"
$ws.Range("M98").Value = 270
$ws.Range("N98").Value = 9500
$ws.Range("O98").Value = 10000
$ws.Range("P98").Value = 9778
$ws.Range("Q98").Value = "`$/bandeja 2 kilos"
$ws.Range("R98").Value = "Provincia de Linares"
$ws.Range("S98").Value = 4889
$ws.Range("T98").Value = 2
